# Auto-generated edit script: updates cryptocurrency price/volume table
# to reflect the GitHub Actions scheduled refresh of cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.240.95"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.871.93"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5020"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09840"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +26.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.123"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.37"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.455"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "1.872.20"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.378"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001135"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.35"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.104"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").Value = "28.334.88"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.259"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.52%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.528"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("D28").Value = "2.087.28"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.94"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.61"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1062"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.631"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.603"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06812"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.386"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02391"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2187"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.011"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.45"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6309"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6006"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.283"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.674"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06848"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.83%  "
